$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wishlist Raluca")

# Add "SKIP" markers to column E for rows 29, 30, 40, 41
$ws.Range("E29").Value = "SKIP"
$ws.Range("E30").Value = "SKIP"
$ws.Range("E40").Value = "SKIP"
$ws.Range("E41").Value = "SKIP"

# Update the active selection on the sheet to E40
$ws.Range("E40").Select()
